$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format date/time-like columns as Text so values such as "2023-09-10"
# are stored as literal strings rather than being auto-converted to Excel
# date/time serial numbers.
$ws.Range("Y21:Y29").NumberFormat = "@"
$ws.Range("Z21:Z29").NumberFormat = "@"
$ws.Range("AA21:AA29").NumberFormat = "@"
$ws.Range("AB21:AB29").NumberFormat = "@"

# Row 21
$ws.Range("A21").Value = 112014423
$ws.Range("B21").Value = 90658
$ws.Range("C21").Value = "Ovaliderad"
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 4361
$ws.Range("F21").Value = "Orange taggsvamp"
$ws.Range("G21").Value = "Hydnellum aurantiacum"
$ws.Range("H21").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P21").Value = "Fläcksberget V, Hjd"
$ws.Range("Q21").Value = 467430.0274016621
$ws.Range("R21").Value = 6875237.811246304
$ws.Range("S21").Value = 20
$ws.Range("T21").Value = "Jämtland"
$ws.Range("U21").Value = "Härjedalen"
$ws.Range("V21").Value = "Härjedalen"
$ws.Range("W21").Value = "Sveg"
$ws.Range("Y21").Value = "2023-09-10"
$ws.Range("Z21").Value = "00:00"
$ws.Range("AA21").Value = "2023-09-10"
$ws.Range("AB21").Value = "00:00"
$ws.Range("AD21").Value = $False
$ws.Range("AE21").Value = $False
$ws.Range("AG21").Value = $False
$ws.Range("AW21").Value = "lennart karlsson"
$ws.Range("AX21").Value = "lennart karlsson"

# Row 22
$ws.Range("A22").Value = 112014208
$ws.Range("B22").Value = 90658
$ws.Range("C22").Value = "Ovaliderad"
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 4361
$ws.Range("F22").Value = "Orange taggsvamp"
$ws.Range("G22").Value = "Hydnellum aurantiacum"
$ws.Range("H22").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P22").Value = "Fläcksberget V, Hjd"
$ws.Range("Q22").Value = 467418.043506761
$ws.Range("R22").Value = 6875312.610613029
$ws.Range("S22").Value = 20
$ws.Range("T22").Value = "Jämtland"
$ws.Range("U22").Value = "Härjedalen"
$ws.Range("V22").Value = "Härjedalen"
$ws.Range("W22").Value = "Sveg"
$ws.Range("Y22").Value = "2023-09-10"
$ws.Range("Z22").Value = "00:00"
$ws.Range("AA22").Value = "2023-09-10"
$ws.Range("AB22").Value = "00:00"
$ws.Range("AD22").Value = $False
$ws.Range("AE22").Value = $False
$ws.Range("AG22").Value = $False
$ws.Range("AW22").Value = "lennart karlsson"
$ws.Range("AX22").Value = "lennart karlsson"

# Row 23
$ws.Range("A23").Value = 112015011
$ws.Range("B23").Value = 90658
$ws.Range("C23").Value = "Ovaliderad"
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 4361
$ws.Range("F23").Value = "Orange taggsvamp"
$ws.Range("G23").Value = "Hydnellum aurantiacum"
$ws.Range("H23").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P23").Value = "Fläcksberget V, Hjd"
$ws.Range("Q23").Value = 467389.9660160011
$ws.Range("R23").Value = 6875327.91063729
$ws.Range("S23").Value = 20
$ws.Range("T23").Value = "Jämtland"
$ws.Range("U23").Value = "Härjedalen"
$ws.Range("V23").Value = "Härjedalen"
$ws.Range("W23").Value = "Sveg"
$ws.Range("Y23").Value = "2023-09-10"
$ws.Range("Z23").Value = "00:00"
$ws.Range("AA23").Value = "2023-09-10"
$ws.Range("AB23").Value = "00:00"
$ws.Range("AD23").Value = $False
$ws.Range("AE23").Value = $False
$ws.Range("AG23").Value = $False
$ws.Range("AW23").Value = "lennart karlsson"
$ws.Range("AX23").Value = "lennart karlsson"

# Row 24
$ws.Range("A24").Value = 112014229
$ws.Range("B24").Value = 90682
$ws.Range("C24").Value = "Ovaliderad"
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 2059
$ws.Range("F24").Value = "Skrovlig taggsvamp"
$ws.Range("G24").Value = "Hydnellum scabrosum"
$ws.Range("H24").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("P24").Value = "Fläcksberget V, Hjd"
$ws.Range("Q24").Value = 467427.230114766
$ws.Range("R24").Value = 6875289.506732536
$ws.Range("S24").Value = 20
$ws.Range("T24").Value = "Jämtland"
$ws.Range("U24").Value = "Härjedalen"
$ws.Range("V24").Value = "Härjedalen"
$ws.Range("W24").Value = "Sveg"
$ws.Range("Y24").Value = "2023-09-10"
$ws.Range("Z24").Value = "00:00"
$ws.Range("AA24").Value = "2023-09-10"
$ws.Range("AB24").Value = "00:00"
$ws.Range("AD24").Value = $False
$ws.Range("AE24").Value = $False
$ws.Range("AG24").Value = $False
$ws.Range("AW24").Value = "lennart karlsson"
$ws.Range("AX24").Value = "lennart karlsson"

# Row 25
$ws.Range("A25").Value = 112014142
$ws.Range("B25").Value = 90666
$ws.Range("C25").Value = "Ovaliderad"
$ws.Range("D25").Value = "LC"
$ws.Range("E25").Value = 4364
$ws.Range("F25").Value = "Dropptaggsvamp"
$ws.Range("G25").Value = "Hydnellum ferrugineum"
$ws.Range("H25").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P25").Value = "Fläcksberget V, Hjd"
$ws.Range("Q25").Value = 467442.7363991642
$ws.Range("R25").Value = 6875336.798642672
$ws.Range("S25").Value = 20
$ws.Range("T25").Value = "Jämtland"
$ws.Range("U25").Value = "Härjedalen"
$ws.Range("V25").Value = "Härjedalen"
$ws.Range("W25").Value = "Sveg"
$ws.Range("Y25").Value = "2023-09-10"
$ws.Range("Z25").Value = "00:00"
$ws.Range("AA25").Value = "2023-09-10"
$ws.Range("AB25").Value = "00:00"
$ws.Range("AD25").Value = $False
$ws.Range("AE25").Value = $False
$ws.Range("AG25").Value = $False
$ws.Range("AW25").Value = "lennart karlsson"
$ws.Range("AX25").Value = "lennart karlsson"

# Row 26
$ws.Range("A26").Value = 112014923
$ws.Range("B26").Value = 90689
$ws.Range("C26").Value = "Ovaliderad"
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 5966
$ws.Range("F26").Value = "Motaggsvamp"
$ws.Range("G26").Value = "Sarcodon squamosus"
$ws.Range("H26").Value = "(Schaeff.) Quél."
$ws.Range("P26").Value = "Fläcksberget V, Hjd"
$ws.Range("Q26").Value = 467413.0579403224
$ws.Range("R26").Value = 6875234.216212902
$ws.Range("S26").Value = 20
$ws.Range("T26").Value = "Jämtland"
$ws.Range("U26").Value = "Härjedalen"
$ws.Range("V26").Value = "Härjedalen"
$ws.Range("W26").Value = "Sveg"
$ws.Range("Y26").Value = "2023-09-10"
$ws.Range("Z26").Value = "00:00"
$ws.Range("AA26").Value = "2023-09-10"
$ws.Range("AB26").Value = "00:00"
$ws.Range("AD26").Value = $False
$ws.Range("AE26").Value = $False
$ws.Range("AG26").Value = $False
$ws.Range("AW26").Value = "lennart karlsson"
$ws.Range("AX26").Value = "lennart karlsson"

# Row 27
$ws.Range("A27").Value = 112014177
$ws.Range("B27").Value = 90689
$ws.Range("C27").Value = "Ovaliderad"
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 5966
$ws.Range("F27").Value = "Motaggsvamp"
$ws.Range("G27").Value = "Sarcodon squamosus"
$ws.Range("H27").Value = "(Schaeff.) Quél."
$ws.Range("P27").Value = "Fläcksberget V, Hjd"
$ws.Range("Q27").Value = 467389.9660160011
$ws.Range("R27").Value = 6875327.91063729
$ws.Range("S27").Value = 20
$ws.Range("T27").Value = "Jämtland"
$ws.Range("U27").Value = "Härjedalen"
$ws.Range("V27").Value = "Härjedalen"
$ws.Range("W27").Value = "Sveg"
$ws.Range("Y27").Value = "2023-09-10"
$ws.Range("Z27").Value = "00:00"
$ws.Range("AA27").Value = "2023-09-10"
$ws.Range("AB27").Value = "00:00"
$ws.Range("AD27").Value = $False
$ws.Range("AE27").Value = $False
$ws.Range("AG27").Value = $False
$ws.Range("AW27").Value = "lennart karlsson"
$ws.Range("AX27").Value = "lennart karlsson"

# Row 28
$ws.Range("A28").Value = 112014300
$ws.Range("B28").Value = 90689
$ws.Range("C28").Value = "Ovaliderad"
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 5966
$ws.Range("F28").Value = "Motaggsvamp"
$ws.Range("G28").Value = "Sarcodon squamosus"
$ws.Range("H28").Value = "(Schaeff.) Quél."
$ws.Range("P28").Value = "Fläcksberget V, Hjd"
$ws.Range("Q28").Value = 467415.4484496959
$ws.Range("R28").Value = 6875287.271149865
$ws.Range("S28").Value = 20
$ws.Range("T28").Value = "Jämtland"
$ws.Range("U28").Value = "Härjedalen"
$ws.Range("V28").Value = "Härjedalen"
$ws.Range("W28").Value = "Sveg"
$ws.Range("Y28").Value = "2023-09-10"
$ws.Range("Z28").Value = "00:00"
$ws.Range("AA28").Value = "2023-09-10"
$ws.Range("AB28").Value = "00:00"
$ws.Range("AD28").Value = $False
$ws.Range("AE28").Value = $False
$ws.Range("AG28").Value = $False
$ws.Range("AW28").Value = "lennart karlsson"
$ws.Range("AX28").Value = "lennart karlsson"

# Row 29
$ws.Range("A29").Value = 112014347
$ws.Range("B29").Value = 90678
$ws.Range("C29").Value = "Ovaliderad"
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 4366
$ws.Range("F29").Value = "Skarp dropptaggsvamp"
$ws.Range("G29").Value = "Hydnellum peckii"
$ws.Range("H29").Value = "Banker"
$ws.Range("P29").Value = "Fläcksberget V, Hjd"
$ws.Range("Q29").Value = 467430.0274016621
$ws.Range("R29").Value = 6875237.811246304
$ws.Range("S29").Value = 20
$ws.Range("T29").Value = "Jämtland"
$ws.Range("U29").Value = "Härjedalen"
$ws.Range("V29").Value = "Härjedalen"
$ws.Range("W29").Value = "Sveg"
$ws.Range("Y29").Value = "2023-09-10"
$ws.Range("Z29").Value = "00:00"
$ws.Range("AA29").Value = "2023-09-10"
$ws.Range("AB29").Value = "00:00"
$ws.Range("AD29").Value = $False
$ws.Range("AE29").Value = $False
$ws.Range("AG29").Value = $False
$ws.Range("AW29").Value = "lennart karlsson"
$ws.Range("AX29").Value = "lennart karlsson"
